# Update countries & provincias Spain
# Applies the 16-Jul-2020 00:25 data refresh to the "Pais" sheet:
#  - updates the "Datos actualizados a ..." timestamp footer
#  - refreshes case/death counters for several countries
#  - re-sorts a handful of rows whose "Casos totales" changed enough to
#    overtake their neighbour (labels + figures swap rows, descending sort
#    by column B is preserved)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Footer timestamp (row 1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 00:25"

# ---------------------------------------------------------------------
# Helper-less, explicit per-row/column writes (row, col, value)
# Columns: A=Pais B=Casos totales C=Nuevos casos D=Casos activos
#          E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes
# ---------------------------------------------------------------------

# Row 4: Estados Unidos - updated counters
$ws.Cells.Item(4, 2).Value = 3608829
$ws.Cells.Item(4, 3).Value = 63752
$ws.Cells.Item(4, 4).Value = 1635913
$ws.Cells.Item(4, 5).Value = 1832938
$ws.Cells.Item(4, 7).Value = 835
$ws.Cells.Item(4, 8).Value = 139978

# Row 5: Brasil - updated counters
$ws.Cells.Item(5, 2).Value = 1966748
$ws.Cells.Item(5, 3).Value = 35544
$ws.Cells.Item(5, 5).Value = 677870
$ws.Cells.Item(5, 7).Value = 1104
$ws.Cells.Item(5, 8).Value = 75366

# Row 30: Suecia - updated counters
$ws.Cells.Item(30, 2).Value = 76492
$ws.Cells.Item(30, 3).Value = 125
$ws.Cells.Item(30, 7).Value = 10
$ws.Cells.Item(30, 8).Value = 5572

# Rows 49-50: Barein overtakes Rumania (Casos totales 34560 > 34226)
# Row 49 becomes Barein with refreshed figures
$ws.Cells.Item(49, 1).Value = "Barein"
$ws.Cells.Item(49, 2).Value = 34560
$ws.Cells.Item(49, 3).Value = 482
$ws.Cells.Item(49, 4).Value = 30320
$ws.Cells.Item(49, 5).Value = 4123
$ws.Cells.Item(49, 7).Value = 6
$ws.Cells.Item(49, 8).Value = 117
# Row 50 becomes Rumania, keeping its previous (unchanged) figures
$ws.Cells.Item(50, 1).Value = "Rumania"
$ws.Cells.Item(50, 2).Value = 34226
$ws.Cells.Item(50, 3).Value = 641
$ws.Cells.Item(50, 4).Value = 22049
$ws.Cells.Item(50, 5).Value = 10225
$ws.Cells.Item(50, 7).Value = 21
$ws.Cells.Item(50, 8).Value = 1952

# Row 67: Uzbekistan - updated counters
$ws.Cells.Item(67, 2).Value = 14581
$ws.Cells.Item(67, 3).Value = 496
$ws.Cells.Item(67, 4).Value = 8655
$ws.Cells.Item(67, 5).Value = 5855

# Rows 75-76: Sudan overtakes Australia (Casos totales 10527 > 10487)
# Row 75 becomes Sudan with refreshed figures
$ws.Cells.Item(75, 1).Value = "Sudan"
$ws.Cells.Item(75, 2).Value = 10527
$ws.Cells.Item(75, 3).Value = 110
$ws.Cells.Item(75, 4).Value = 5601
$ws.Cells.Item(75, 5).Value = 4258
$ws.Cells.Item(75, 7).Value = 9
$ws.Cells.Item(75, 8).Value = 668
# Row 76 becomes Australia, keeping its previous (unchanged) figures
$ws.Cells.Item(76, 1).Value = "Australia"
$ws.Cells.Item(76, 2).Value = 10487
$ws.Cells.Item(76, 3).Value = 237
$ws.Cells.Item(76, 4).Value = 7928
$ws.Cells.Item(76, 5).Value = 2448
$ws.Cells.Item(76, 7).Value = 3
$ws.Cells.Item(76, 8).Value = 111

# Row 85: Bulgaria - updated counters
$ws.Cells.Item(85, 2).Value = 7877
$ws.Cells.Item(85, 3).Value = 232
$ws.Cells.Item(85, 4).Value = 3841
$ws.Cells.Item(85, 5).Value = 3747
$ws.Cells.Item(85, 7).Value = 6
$ws.Cells.Item(85, 8).Value = 289

# Row 91: Guayana Francesa - updated counters
$ws.Cells.Item(91, 2).Value = 6299
$ws.Cells.Item(91, 3).Value = 70
$ws.Cells.Item(91, 4).Value = 3738
$ws.Cells.Item(91, 5).Value = 2528
$ws.Cells.Item(91, 7).Value = 2
$ws.Cells.Item(91, 8).Value = 33

# Rows 104-106: Paraguay jumps above Nicaragua and Somalia
# (Casos totales 3198 > 3147 > 3076)
# Row 104 becomes Paraguay with refreshed figures
$ws.Cells.Item(104, 1).Value = "Paraguay"
$ws.Cells.Item(104, 2).Value = 3198
$ws.Cells.Item(104, 3).Value = 124
$ws.Cells.Item(104, 4).Value = 1338
$ws.Cells.Item(104, 5).Value = 1835
$ws.Cells.Item(104, 8).Value = 25
# Row 105 becomes Nicaragua, keeping its previous (unchanged) figures
$ws.Cells.Item(105, 1).Value = "Nicaragua"
$ws.Cells.Item(105, 2).Value = 3147
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = 2282
$ws.Cells.Item(105, 5).Value = 766
$ws.Cells.Item(105, 8).Value = 99
# Row 106 becomes Somalia, keeping its previous (unchanged) figures
$ws.Cells.Item(106, 1).Value = "Somalia"
$ws.Cells.Item(106, 2).Value = 3076
$ws.Cells.Item(106, 3).Value = 0
$ws.Cells.Item(106, 4).Value = 1380
$ws.Cells.Item(106, 5).Value = 1603
$ws.Cells.Item(106, 8).Value = 93

# Rows 111-112: Malaui overtakes Libano (Casos totales 2614 > 2542)
# Row 111 becomes Malaui with refreshed figures
$ws.Cells.Item(111, 1).Value = "Malaui"
$ws.Cells.Item(111, 2).Value = 2614
$ws.Cells.Item(111, 3).Value = 117
$ws.Cells.Item(111, 4).Value = 1005
$ws.Cells.Item(111, 5).Value = 1566
$ws.Cells.Item(111, 7).Value = 3
$ws.Cells.Item(111, 8).Value = 43
# Row 112 becomes Libano, keeping its previous (unchanged) figures
$ws.Cells.Item(112, 1).Value = "Libano"
$ws.Cells.Item(112, 2).Value = 2542
$ws.Cells.Item(112, 3).Value = 91
$ws.Cells.Item(112, 4).Value = 1455
$ws.Cells.Item(112, 5).Value = 1049
$ws.Cells.Item(112, 7).Value = 1
$ws.Cells.Item(112, 8).Value = 38

# Rows 209-210: Groenlandia / Islas Malvinas tie (both 13) swap order
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"
